$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row with the new Hardware-ID and its expiry date
$ws.Range("A3").Value = "S36SNWAH859775X"

# Serial date value for 2026-02-11 (matches existing date column's format)
$ws.Range("B3").Value = 46064

# Match the date formatting used by the existing date cell (B2)
$ws.Range("B2").Copy()
$ws.Range("B3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update selection to reflect the new active cell (next empty row)
$ws.Range("B4").Select()
